$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix rendering issues with Fall 2023 update

# 1. Add missing space after <br> in the D2 reading text
$ws.Range("D2").Value = "LSWR Ch 2 and 3 <br> Light, Singer & Willet 1990, Ch. 2"

# 2. "Units" -> "Unit" in the header cell A1
$ws.Range("A1").Value = "Unit"

# 3. Update the active selection from D3 to A2
$ws.Range("A2").Select()
